# Sync automático del tracker (cada 3h)
# Appends the latest batch of matches/pronósticos to the bottom of the
# tracker sheet (Sheet1), starting right after the last existing row (413).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ EventId = "14528379"; Fecha = "2025-08-25"; JugadorA = "Gonzalo Bueno";             JugadorB = "Facundo Diaz Acosta";          Pronostico = "Gana Gonzalo Bueno";                Cuota = 3.25 },
    @{ EventId = "14528375"; Fecha = "2025-08-25"; JugadorA = "Jacopo Berrettini";          JugadorB = "Luciano Emanuel Ambrogi";       Pronostico = "Gana Luciano Emanuel Ambrogi";      Cuota = 2.38 },
    @{ EventId = "14527264"; Fecha = "2025-08-25"; JugadorA = "Jie Cui";                    JugadorB = "Arthur Gea";                    Pronostico = "Gana Jie Cui";                      Cuota = 3.75 },
    @{ EventId = "14528390"; Fecha = "2025-08-25"; JugadorA = "Alex Barrena";               JugadorB = "Christoph Negritu";             Pronostico = "Gana Christoph Negritu";            Cuota = 2.38 },
    @{ EventId = "14528389"; Fecha = "2025-08-25"; JugadorA = "Nicolas Alvarez Varona";     JugadorB = "Robert Strombachs";             Pronostico = "Gana Robert Strombachs";            Cuota = 2.63 },
    @{ EventId = "14528387"; Fecha = "2025-08-25"; JugadorA = "Pedro Araujo";               JugadorB = "Elias Ymer";                    Pronostico = "Gana Pedro Araujo";                 Cuota = 6 },
    @{ EventId = "14528393"; Fecha = "2025-08-25"; JugadorA = "Thiago Seyboth Wild";        JugadorB = "Nikolas Sanchez Izquierdo";     Pronostico = "Gana Nikolas Sanchez Izquierdo";    Cuota = 2.75 }
)

$startRow = 414
$endRow = $startRow + $newRows.Count - 1

# event_id (A) and fecha (B) travel as plain text in the source feed (the
# id is too big/arbitrary to treat as a number and the date must stay as
# "yyyy-mm-dd" instead of becoming a serial date), so mark that block as
# Text before writing into it.
$ws.Range("A$startRow`:B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.EventId
    $ws.Cells.Item($r, 2).Value = $row.Fecha
    $ws.Cells.Item($r, 3).Value = $row.JugadorA
    $ws.Cells.Item($r, 4).Value = $row.JugadorB
    $ws.Cells.Item($r, 5).Value = $row.Pronostico
    $ws.Cells.Item($r, 6).Value = $row.Cuota
    # resultado (G) / profit (H) are left blank until the match is settled.
}
